$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.577.18"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.878.59"
$ws.Range("E3").Value = "  -0.10%  "

# Row 5 - BNB
$ws.Range("D5").Value = "249.48"
$ws.Range("E5").Value = "  +1.82%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.04%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4751"
$ws.Range("E7").Value = "  -0.22%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2937"
$ws.Range("E8").Value = "  +1.41%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06529"
$ws.Range("E9").Value = "  +0.13%  "

# Row 10 - Solana
$ws.Range("D10").Value = "22.15"
$ws.Range("E10").Value = "  +3.48%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.07735"
$ws.Range("E11").Value = "  +0.22%  "

# Row 12 - Polygon
$ws.Range("D12").Value = "0.7407"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13 - Litecoin
$ws.Range("D13").Value = "96.81"
$ws.Range("E13").Value = "  -0.45%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.877.38"
$ws.Range("E14").Value = "  +0.12%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "5.245"
$ws.Range("E15").Value = "  +2.22%  "

# Row 16 - BitcoinCash
$ws.Range("D16").Value = "273.41"
$ws.Range("E16").Value = "  -0.56%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "30.716.69"
$ws.Range("E17").Value = "  +0.95%  "

# Row 18 - Avalanche
$ws.Range("D18").Value = "13.24"
$ws.Range("E18").Value = "  -2.11%  "

# Row 19 - Dai
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.000007529"
$ws.Range("E20").Value = "  -0.08%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.124.50"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22 - was BinanceUSD, now Uniswap
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "5.290"
$ws.Range("E22").Value = "  +0.75%  "

# Row 23 - was Uniswap, now BinanceUSD
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "6.217"
$ws.Range("E24").Value = "  +0.77%  "

# Row 25 - was Cosmos, now Monero
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "165.09"
$ws.Range("E25").Value = "  +0.58%  "

# Row 26 - was Monero, now Cosmos
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "9.181"
$ws.Range("E26").Value = "  -0.54%  "

# Row 27 - EthereumClassic
$ws.Range("E27").Value = "  -0.28%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "1.912"
$ws.Range("E28").Value = "  -2.14%  "

# Row 29 - Stellar
$ws.Range("D29").Value = "0.09797"
$ws.Range("E29").Value = "  -1.40%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  -2.23%  "

# Row 31 - PancakeSwap
$ws.Range("D31").Value = "1.505"
$ws.Range("E31").Value = "  -0.30%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.288"
$ws.Range("E32").Value = "  -0.68%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "4.126"
$ws.Range("E33").Value = "  +1.45%  "

# Row 34 - Hedera
$ws.Range("D34").Value = "0.04891"
$ws.Range("E34").Value = "  +3.71%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.126"
$ws.Range("E35").Value = "  +0.39%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.6960"
$ws.Range("E36").Value = "  -0.17%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value = "2.720"

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01890"
$ws.Range("E38").Value = "  +1.99%  "

# Row 39 - MXToken
$ws.Range("D39").Value = "2.766"
$ws.Range("E39").Value = "  +0.17%  "

# Row 40 - FraxShare
$ws.Range("E40").Value = "  +0.00%  "

# Row 41 - Aave
$ws.Range("D41").Value = "74.11"
$ws.Range("E41").Value = "  +6.78%  "

# Row 42 - RenderToken
$ws.Range("D42").Value = "2.007"
$ws.Range("E42").Value = "  +5.30%  "

# Row 43 - TheSandbox
$ws.Range("D43").Value = "0.4256"
$ws.Range("E43").Value = "  +2.00%  "

# Row 44 - PaxDollar
$ws.Range("E44").Value = "  +0.06%  "

# Row 45 - TrustWalletToken
$ws.Range("D45").Value = "0.8393"
$ws.Range("E45").Value = "  -0.05%  "

# Row 46 - Quant
$ws.Range("D46").Value = "102.20"
$ws.Range("E46").Value = "  +0.10%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "9.418"
$ws.Range("E47").Value = "  +2.44%  "

# Row 48 - Aptos
$ws.Range("D48").Value = "7.055"
$ws.Range("E48").Value = "  -0.48%  "

# Row 49 - Elrond
$ws.Range("D49").Value = "35.48"
$ws.Range("E49").Value = "  +0.71%  "

# Row 50 - Maker
$ws.Range("D50").Value = "917.56"
$ws.Range("E50").Value = "  -0.63%  "

# Row 51 - was Decentraland, now Cronos
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05716"
$ws.Range("E51").Value = "  +2.23%  "
